$wb = $excel.ActiveWorkbook

# Sheet ALC, row 2 (Leve Item ID 5489)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 205.36363
$ws.Range("I2").Value = 250
$ws.Range("J2").Value = 151.8
$ws.Range("K2").Value = 250
$ws.Range("L2").Value = 151.8
$ws.Range("M2").Value = -137
$ws.Range("N2").Value = -377.8

# Sheet ALC, row 5 (Leve Item ID 5503)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 95.86667
$ws.Range("I5").Value = 72
$ws.Range("J5").Value = 131.66667
$ws.Range("K5").Value = 72
$ws.Range("L5").Value = 131.66667
$ws.Range("M5").Value = 43
$ws.Range("N5").Value = -361.66667

# Sheet ALC, row 76 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3128.2068
$ws.Range("I76").Value = 3160.72
$ws.Range("J76").Value = 2925
$ws.Range("K76").Value = 3160.72
$ws.Range("L76").Value = 2925
$ws.Range("M76").Value = -2845.72
$ws.Range("N76").Value = -3555

# Sheet ALC, row 79 (Leve Item ID 12602)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3128.2068
$ws.Range("I79").Value = 3160.72
$ws.Range("J79").Value = 2925
$ws.Range("K79").Value = 3160.72
$ws.Range("L79").Value = 2925
$ws.Range("M79").Value = -2068.72
$ws.Range("N79").Value = -5109

# Sheet ALC, row 86 (Leve Item ID 12603)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2800.8965
$ws.Range("I86").Value = 1724.8096
$ws.Range("J86").Value = 5625.625
$ws.Range("K86").Value = 1724.8096
$ws.Range("L86").Value = 5625.625
$ws.Range("M86").Value = -601.8096
$ws.Range("N86").Value = -7871.625

# Sheet ALC, row 89 (Leve Item ID 12603)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 2800.8965
$ws.Range("I89").Value = 1724.8096
$ws.Range("J89").Value = 5625.625
$ws.Range("K89").Value = 8624.048000000001
$ws.Range("L89").Value = 28128.125
$ws.Range("M89").Value = -3008.048000000001
$ws.Range("N89").Value = -39360.125

# Sheet ALC, row 96 (Leve Item ID 19894)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 13274.75
$ws.Range("I96").Value = 686.6667
$ws.Range("J96").Value = 20827.6
$ws.Range("K96").Value = 2060.0001
$ws.Range("L96").Value = 62482.8
$ws.Range("M96").Value = -687.0001000000002
$ws.Range("N96").Value = -65228.8

# Sheet ALC, row 100 (Leve Item ID 19906)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1509.5238
$ws.Range("I100").Value = 1367.4286
$ws.Range("J100").Value = 1793.7142
$ws.Range("K100").Value = 1367.4286
$ws.Range("L100").Value = 1793.7142
$ws.Range("M100").Value = -826.4286
$ws.Range("N100").Value = -2875.7142

# Sheet ALC, row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 265830.1
$ws.Range("I132").Value = 296989.5
$ws.Range("J132").Value = 975
$ws.Range("K132").Value = 890968.5
$ws.Range("L132").Value = 2925
$ws.Range("M132").Value = -888438.5
$ws.Range("N132").Value = -7985

# Sheet ALC, row 137 (Leve Item ID 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 26317382
$ws.Range("I137").Value = 1056.8
$ws.Range("J137").Value = 125003600
$ws.Range("K137").Value = 3170.4
$ws.Range("L137").Value = 375010800
$ws.Range("M137").Value = -620.3999999999996
$ws.Range("N137").Value = -375015900

# Sheet ARM, row 2 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4203205.5
$ws.Range("J2").Value = 5884091.5
$ws.Range("L2").Value = 5884091.5
$ws.Range("N2").Value = -5884317.5

# Sheet ARM, row 32 (Leve Item ID 44147)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4695.557
$ws.Range("I32").Value = 4771.4653
$ws.Range("J32").Value = 4328.6665
$ws.Range("K32").Value = 4771.4653
$ws.Range("L32").Value = 4328.6665
$ws.Range("M32").Value = -4484.4653
$ws.Range("N32").Value = -4902.6665

# Sheet ARM, row 116 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 4203205.5
$ws.Range("J116").Value = 5884091.5
$ws.Range("L116").Value = 5884091.5
$ws.Range("N116").Value = -5888679.5

# Sheet ARM, row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1239.4038
$ws.Range("I132").Value = 1124.5641
$ws.Range("J132").Value = 1583.9231
$ws.Range("K132").Value = 3373.6923
$ws.Range("L132").Value = 4751.7693
$ws.Range("M132").Value = -843.6923000000002
$ws.Range("N132").Value = -9811.7693

# Sheet BSM, row 3 (Leve Item ID 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4203205.5
$ws.Range("J3").Value = 5884091.5
$ws.Range("L3").Value = 5884091.5
$ws.Range("N3").Value = -5884319.5

# Sheet BSM, row 86 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1991.7778
$ws.Range("I86").Value = 1875.5294
$ws.Range("J86").Value = 2189.4
$ws.Range("K86").Value = 1875.5294
$ws.Range("L86").Value = 2189.4
$ws.Range("M86").Value = -752.5293999999999
$ws.Range("N86").Value = -4435.4

# Sheet BSM, row 89 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1991.7778
$ws.Range("I89").Value = 1875.5294
$ws.Range("J89").Value = 2189.4
$ws.Range("K89").Value = 9377.646999999999
$ws.Range("L89").Value = 10947
$ws.Range("M89").Value = -3761.646999999999
$ws.Range("N89").Value = -22179

# Sheet BSM, row 105 (Leve Item ID 19947)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5128.5713
$ws.Range("I105").Value = 5128.5713
$ws.Range("K105").Value = 5128.5713
$ws.Range("M105").Value = -3381.5713

# Sheet CRP, row 62 (Leve Item ID 12580)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2688.75
$ws.Range("I62").Value = 2402.5
$ws.Range("J62").Value = 2975
$ws.Range("K62").Value = 2402.5
$ws.Range("L62").Value = 2975
$ws.Range("M62").Value = -1778.5
$ws.Range("N62").Value = -4223

# Sheet CRP, row 65 (Leve Item ID 12580)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2688.75
$ws.Range("I65").Value = 2402.5
$ws.Range("J65").Value = 2975
$ws.Range("K65").Value = 12012.5
$ws.Range("L65").Value = 14875
$ws.Range("M65").Value = -8892.5
$ws.Range("N65").Value = -21115

# Sheet CUL, row 2 (Leve Item ID 4847)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 582775.5
$ws.Range("I2").Value = 30.533333
$ws.Range("J2").Value = 1831514.8
$ws.Range("K2").Value = 183.199998
$ws.Range("L2").Value = 10989088.8
$ws.Range("M2").Value = -70.19999799999999
$ws.Range("N2").Value = -10989314.8

# Sheet CUL, row 9 (Leve Item ID 4681)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1233
$ws.Range("I9").Value = 830
$ws.Range("J9").Value = 1333.75
$ws.Range("K9").Value = 2490
$ws.Range("L9").Value = 4001.25
$ws.Range("M9").Value = -2266
$ws.Range("N9").Value = -4449.25

# Sheet CUL, row 13 (Leve Item ID 4657)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 68
$ws.Range("I13").Value = 75
$ws.Range("J13").Value = 40
$ws.Range("K13").Value = 225
$ws.Range("L13").Value = 120
$ws.Range("M13").Value = -57
$ws.Range("N13").Value = -456

# Sheet CUL, row 40 (Leve Item ID 4827)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 513.2727
$ws.Range("I40").Value = 72.90909000000001
$ws.Range("J40").Value = 953.63635
$ws.Range("K40").Value = 291.63636
$ws.Range("L40").Value = 3814.5454
$ws.Range("M40").Value = -222.63636
$ws.Range("N40").Value = -3952.5454

# Sheet CUL, row 104 (Leve Item ID 19807)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 3299
$ws.Range("J104").Value = 3299
$ws.Range("L104").Value = 9897
$ws.Range("N104").Value = -15139

# Sheet CUL, row 122 (Leve Item ID 36078)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 23810860
$ws.Range("I122").Value = 37037396
$ws.Range("K122").Value = 333336564
$ws.Range("M122").Value = -333334114

# Sheet CUL, row 126 (Leve Item ID 36045)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 8676.137000000001
$ws.Range("I126").Value = 963.3333
$ws.Range("J126").Value = 9893.947
$ws.Range("K126").Value = 2889.9999
$ws.Range("L126").Value = 29681.841
$ws.Range("M126").Value = 2050.0001
$ws.Range("N126").Value = -39561.841

# Sheet CUL, row 129 (Leve Item ID 36054)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1690.4286
$ws.Range("I129").Value = 400
$ws.Range("J129").Value = 1905.5
$ws.Range("K129").Value = 1200
$ws.Range("L129").Value = 5716.5
$ws.Range("M129").Value = 3800
$ws.Range("N129").Value = -15716.5

# Sheet CUL, row 131 (Leve Item ID 36060)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1206792.6
$ws.Range("J131").Value = 1472191
$ws.Range("L131").Value = 4416573
$ws.Range("N131").Value = -4426653

# Sheet GSM, row 80 (Leve Item ID 12521)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 51500
$ws.Range("J80").Value = 51500
$ws.Range("L80").Value = 51500
$ws.Range("N80").Value = -53496

# Sheet GSM, row 83 (Leve Item ID 12521)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 51500
$ws.Range("J83").Value = 51500
$ws.Range("L83").Value = 257500
$ws.Range("N83").Value = -267484

# Sheet GSM, row 113 (Leve Item ID 27710)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1290.0769
$ws.Range("I113").Value = 763.55554
$ws.Range("J113").Value = 2474.75
$ws.Range("K113").Value = 763.55554
$ws.Range("L113").Value = 2474.75
$ws.Range("M113").Value = 1406.44446
$ws.Range("N113").Value = -6814.75

# Sheet LTW, row 2 (Leve Item ID 2631)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 15004001
$ws.Range("J2").Value = 15004001
$ws.Range("L2").Value = 15004001
$ws.Range("N2").Value = -15004225

# Sheet LTW, row 7 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1279.7142
$ws.Range("I7").Value = 981
$ws.Range("K7").Value = 981
$ws.Range("M7").Value = -869

# Sheet LTW, row 36 (Leve Item ID 34261)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0  # special-cased, see analysis
$ws.Range("N36").ClearContents()

# Sheet LTW, row 126 (Leve Item ID 36249)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1279.7142
$ws.Range("I126").Value = 981
$ws.Range("K126").Value = 2943
$ws.Range("M126").Value = -473

# Sheet WVR, row 81 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2922
$ws.Range("I81").Value = 1085.7142
$ws.Range("K81").Value = 2171.4284
$ws.Range("M81").Value = -1110.4284

# Sheet WVR, row 84 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2922
$ws.Range("I84").Value = 1085.7142
$ws.Range("K84").Value = 10857.142
$ws.Range("M84").Value = -5553.142
